# Reklamacje workbook update:
#  - Reorder customer-info columns E:J into (Imie, Nazwisko, Telefon, Email, Indeksy, Numery FedEx)
#  - Rename headers for M/N, shorten the M/N date-time strings to date-only
#  - Update row 4 status + completion date
#  - Append new complaint rows 5-8
#
# NOTE: plain `$cell.Value = "123"` auto-coerces purely-numeric-looking
# strings (phone numbers, tracking numbers, IDs, date strings, ...) into
# Number/Date cells, which loses leading "+"/zeros and reformats dates.
# Forcing NumberFormat to "@" (Text) before the assignment keeps the
# underlying cell type as Text; resetting the Style back to "Normal"
# afterwards drops the now-unneeded text format so the cell keeps the
# default (unstyled) look, matching plain data cells elsewhere in the sheet.
function Set-TextValue($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) - rename / reorder columns E..J, and M/N
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 5).Value  = "Imię"          # E1
$ws.Cells.Item(1, 6).Value  = "Nazwisko"      # F1
$ws.Cells.Item(1, 7).Value  = "Telefon"       # G1
$ws.Cells.Item(1, 8).Value  = "Email"         # H1
$ws.Cells.Item(1, 9).Value  = "Indeksy"       # I1
$ws.Cells.Item(1, 10).Value = "Numery FedEx"  # J1
$ws.Cells.Item(1, 13).Value = "Zgłoszono"     # M1
$ws.Cells.Item(1, 14).Value = "Zrealizowano"  # N1

# ---------------------------------------------------------------------------
# 2. Existing data rows 2-4: shuffle E..J into their new positions and trim
#    the timestamp strings in M/N down to the date portion only.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 4; $r++) {
    $oldIndeksy  = $ws.Cells.Item($r, 5).Text   # was E (Indeksy produktow)
    $oldFedEx    = $ws.Cells.Item($r, 6).Text   # was F (Numery FedEx)
    $oldImie     = $ws.Cells.Item($r, 7).Text   # was G (Imie klienta)
    $oldNazwisko = $ws.Cells.Item($r, 8).Text   # was H (Nazwisko klienta)
    $oldTelefon  = $ws.Cells.Item($r, 9).Text   # was I (Telefon klienta)
    $oldEmail    = $ws.Cells.Item($r, 10).Text  # was J (Email klienta)

    Set-TextValue $ws $r 5  $oldImie
    Set-TextValue $ws $r 6  $oldNazwisko
    Set-TextValue $ws $r 7  $oldTelefon
    Set-TextValue $ws $r 8  $oldEmail
    Set-TextValue $ws $r 9  $oldIndeksy
    Set-TextValue $ws $r 10 $oldFedEx

    $oldCreated = $ws.Cells.Item($r, 13).Text
    if ($oldCreated) {
        $createdParts = $oldCreated.Split(" ")
        $createdDate = $createdParts[0]
        Set-TextValue $ws $r 13 $createdDate
    }
    $oldDone = $ws.Cells.Item($r, 14).Text
    if ($oldDone) {
        $doneParts = $oldDone.Split(" ")
        $doneDate = $doneParts[0]
        Set-TextValue $ws $r 14 $doneDate
    }
}

# Row 4 also changes status to completed and now has a completion date
$ws.Cells.Item(4, 12).Value = "Zrealizowana"
Set-TextValue $ws 4 14 "2025-03-20"

# ---------------------------------------------------------------------------
# 3. Append new complaint rows 5-8
# ---------------------------------------------------------------------------
$newRows = @(
    @{ A=5; B="Szafka Nako 3S Biały"; C="Uszkodzone 2 elementy."; D="18311";
       E="Anna "; F="Kostecka"; G="+48577849782"; H="zn09fg1yop+1efe357e8@allegromail.pl";
       I="EL-NAK-3SZ-BOK-P-(F18-)-BI, EL-NAK-3SZ-FRO-SZ-(F33-)-BI"; J=", ";
       K="LP: 6233419925068, 6233423003608"; L="Zrealizowana"; M="2025-03-20"; N="2025-03-20" },

    @{ A=6; B="Komoda Gabriel 6S"; C="Klient nie otrzymał 6 boczków szuflady w swojej paczce."; D="18246";
       E="Ola"; F="Szpatowska"; G="+48691137854"; H="alekssbq@gmail.com";
       I="EL-GAB-BOCZ-SZ-P-BISM, EL-GAB-BOCZ-SZ-P-BISM, EL-GAB-BOCZ-SZ-P-BISM, EL-GAB-BOCZ-SZ-P-BISM, EL-GAB-BOCZ-SZ-P-BISM, EL-GAB-BOCZ-SZ-P-BISM"; J=$null;
       K="6233419673396"; L="Zrealizowana"; M="2025-03-20"; N="2025-03-20" },

    @{ A=7; B="Komoda Gabriel 6S"; C="Uszkodzony Lewy front szuflady."; D="17873";
       E="Anna "; F="Gajko"; G="+48696831228"; H="ania23923@gmail.com";
       I="EL-GAB-6S-FRONT-SZ-L-BISM"; J=$null;
       K=$null; L="Zrealizowana"; M="2025-03-20"; N="2025-03-20" },

    @{ A=8; B="Szafka Nako 3S Biały"; C="Uszkodzone 2 elementy"; D="18800";
       E="Svitlana "; F="Zdybel"; G="+48888036037"; H="5a8tg46sy3+717cd2316@allegromail.pl";
       I="EL-NAK-3SZ-BOK-L-(F18-)-BI, EL-NAK-3SZ-BOK-P-(F18-)-BI"; J=$null;
       K=$null; L="W trakcie"; M="2025-03-20"; N=$null }
)

foreach ($row in $newRows) {
    $r = $row.A
    $ws.Cells.Item($r, 1).Value  = $row.A          # numeric ID
    $ws.Cells.Item($r, 2).Value  = $row.B           # B - plain text, no digits-only risk
    $ws.Cells.Item($r, 3).Value  = $row.C           # C - plain text
    Set-TextValue  $ws $r 4  $row.D                 # D - numeric-looking order number -> force text
    Set-TextValue  $ws $r 5  $row.E                 # E - Imie
    Set-TextValue  $ws $r 6  $row.F                 # F - Nazwisko
    Set-TextValue  $ws $r 7  $row.G                 # G - Telefon (numeric-looking) -> force text
    Set-TextValue  $ws $r 8  $row.H                 # H - Email
    Set-TextValue  $ws $r 9  $row.I                 # I - Indeksy
    if ($row.J) { Set-TextValue $ws $r 10 $row.J }  # J - Numery FedEx (numeric-looking)
    if ($row.K) { Set-TextValue $ws $r 11 $row.K }  # K - Notatki (numeric-looking in some rows)
    $ws.Cells.Item($r, 12).Value = $row.L            # L - Status
    Set-TextValue  $ws $r 13 $row.M                  # M - Zgloszono (date-looking)
    if ($row.N) { Set-TextValue $ws $r 14 $row.N }   # N - Zrealizowano (date-looking)
}
